$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had a bogus first row that repeated the dataset's
# long title string in every column (A1:F1). The real header row (years
# label + category names) lived in row 2, and the yearly data started at
# row 3. Deleting row 1 promotes row 2 to be the new header row 1 and
# shifts every data row up by one, matching the new A1:F21 layout.
$ws.Rows(1).Delete()

# Re-point the sheet view's selection at the new header row instead of
# the old K11 active cell.
$ws.Range("A1:XFD1").Select()

# Make the new header row's style match the target: bold, 11pt (was the
# old 14pt "title" font), black text, no fill.
$ws.Range("B1:F1").Font.Bold = $true
$ws.Range("B1:F1").Font.Size = 11
$ws.Range("B1:F1").Font.Color = 0

# The leading "year" cells in column A use the same bold header style.
$ws.Range("A2:A21").Font.Bold = $true
$ws.Range("A2:A21").Font.Size = 11
$ws.Range("A2:A21").Font.Color = 0

# The ".." placeholder cells (missing data) keep their orange fill and
# right alignment, but drop the bold weight and switch to the theme text
# color instead of explicit black.
$ws.Range("F2:F8").Font.Bold = $false
$ws.Range("F2:F8").Font.ThemeColor = 1
$ws.Range("F2:F8").Interior.Color = 8036607
$ws.Range("F2:F8").HorizontalAlignment = -4152

# Printer/page setup: portrait orientation.
$ws.PageSetup.Orientation = 1

$wb.Save()
